$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header strings (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Update column C data values (imputed GDP figures) ---
$ws.Range("C2").Value = 2934.187009790061
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 1460.056109840828
$ws.Range("C6").Value = 7772.38875590225
$ws.Range("C7").Value = 1909.084588129339
$ws.Range("C8").Value = 6128.19547247793
$ws.Range("C9").Value = 3972.630273980753
$ws.Range("C10").Value = 4729.735976516416
$ws.Range("C11").Value = 471.9591970298227
$ws.Range("C12").Value = 471.181692645893
$ws.Range("C14").Value = 743.403784726004
$ws.Range("C15").Value = 1489.459305835934
$ws.Range("C16").Value = 1286.515571617672
$ws.Range("C17").Value = 892.5687203369533
$ws.Range("C18").Value = 1036.533951644687
$ws.Range("C19").Value = 1213.112645064426
$ws.Range("C20").Value = 1299.344949460393
$ws.Range("C21").Value = 513.4456986202239
$ws.Range("C22").Value = 2292.445156190819
$ws.Range("C23").Value = 951.6879611168786
$ws.Range("C24").Value = 612.3436990512633
$ws.Range("C25").Value = 1401.47747416771
$ws.Range("C26").Value = 2983.242707849043
$ws.Range("C27").Value = 2898.942214704482
$ws.Range("C28").Value = 665.6274194933962
$ws.Range("C29").Value = 1904.346464968814
$ws.Range("C30").Value = 1955.461557360978
$ws.Range("C31").Value = 4633.590358399045
$ws.Range("C32").Value = 5082.354756663512
$ws.Range("C33").Value = 13113.52596172773
$ws.Range("C34").Value = 466.3491752969743
$ws.Range("C35").Value = 1525.562493537689
$ws.Range("C36").Value = 1303.425880277445
$ws.Range("C37").Value = 863.7612548677739
$ws.Range("C38").Value = 2828.483778716848
$ws.Range("C39").Value = 505.2384587280311
$ws.Range("C40").Value = 1037.747039954749
$ws.Range("C41").Value = 1132.548400540401
$ws.Range("C42").Value = 1446.371630707023
$ws.Range("C43").Value = 982.980837581714
$ws.Range("C44").Value = 1591.56825353313
$ws.Range("C45").Value = 864.5379000312432
$ws.Range("C46").Value = 3083.80337578809
$ws.Range("C47").Value = 2965.153206179127
$ws.Range("C48").Value = 691.8942672110555
$ws.Range("C49").Value = 1939.33862702996
$ws.Range("C50").Value = 1577.487171555845
$ws.Range("C51").Value = 2024.117324382548
$ws.Range("C52").Value = 6711.616186806423
$ws.Range("C53").Value = 4921.848409120176
$ws.Range("C54").Value = 4479.398934239905
$ws.Range("C55").Value = 5360.226632400601
$ws.Range("C56").Value = 1410.426304742003
$ws.Range("C57").Value = 13389.959626871
$ws.Range("C58").Value = 8415.999185225515
$ws.Range("C59").Value = 1057.667740311969
$ws.Range("C60").Value = 528.6449273841434
$ws.Range("C61").Value = 909.3123437708064
$ws.Range("C62").Value = 1223.631935023299
$ws.Range("C63").Value = 1543.763984230257
$ws.Range("C64").Value = 564.5208442217756
$ws.Range("C65").Value = 1745.10167474004
$ws.Range("C66").Value = 1000.829216794104
$ws.Range("C67").Value = 788.439151581443
$ws.Range("C68").Value = 869.6014949562591
$ws.Range("C69").Value = 1591.319557098113
$ws.Range("C70").Value = 2094.024217383061
$ws.Range("C71").Value = 5122.180090208862
$ws.Range("C72").Value = 3156.723844635973
$ws.Range("C73").Value = 2999.422762626143
$ws.Range("C74").Value = 1982.009737844954
$ws.Range("C75").Value = 4394.543881413723
$ws.Range("C76").Value = 5642.578115155247
$ws.Range("C77").Value = 865.7498910537106
$ws.Range("C78").Value = 2379.668184479739
$ws.Range("C79").Value = 0
$ws.Range("C80").Value = 13261.19696291444
$ws.Range("C81").Value = 1102.527430026863
$ws.Range("C82").Value = 335.38915520098
$ws.Range("C83").Value = 929.4690557368662
$ws.Range("C84").Value = 1299.811672673934
$ws.Range("C85").Value = 1618.597849849475
$ws.Range("C86").Value = 597.3813896804552
$ws.Range("C87").Value = 1778.60982580794
$ws.Range("C88").Value = 1032.277326842402
$ws.Range("C89").Value = 817.1226340535979
$ws.Range("C90").Value = 872.1235974568563
$ws.Range("C91").Value = 1620.124515672545
$ws.Range("C92").Value = 2201.396847776877
$ws.Range("C93").Value = 5295.682695961288
$ws.Range("C94").Value = 3212.740625904757
$ws.Range("C95").Value = 3056.152683606517
$ws.Range("C96").Value = 2000.792448761861
$ws.Range("C97").Value = 4699.493713911862
$ws.Range("C98").Value = 5919.20956823756
$ws.Range("C99").Value = 907.2574180443885
$ws.Range("C100").Value = 2497.68592515536
$ws.Range("C101").Value = 0
$ws.Range("C102").Value = 13558.79747639988
$ws.Range("C103").Value = 951.3148210424945
$ws.Range("C104").Value = 1140.447753778042
$ws.Range("C105").Value = 2286.013198234259
$ws.Range("C106").Value = 1401.753174264641
$ws.Range("C107").Value = 961.3778847738438
$ws.Range("C108").Value = 1379.14068216006
$ws.Range("C109").Value = 5412.131646018807
$ws.Range("C110").Value = 3252.634165082374
$ws.Range("C111").Value = 1627.760281433693
$ws.Range("C112").Value = 3137.260298393558
$ws.Range("C113").Value = 2025.814194788851
$ws.Range("C114").Value = 1640.18070024053
$ws.Range("C115").Value = 1060.095015975378
$ws.Range("C116").Value = 6753.607115829548
$ws.Range("C117").Value = 586.2293607842975
$ws.Range("C118").Value = 468.1130345750273
$ws.Range("C119").Value = 558.2093442539386
$ws.Range("C120").Value = 711.3043470146426
$ws.Range("C121").Value = 4861.287098802361
$ws.Range("C122").Value = 5996.49696468919
$ws.Range("C123").Value = 0
$ws.Range("C124").Value = 7582.696928894958
$ws.Range("C125").Value = 846.386841468855
$ws.Range("C126").Value = 1644.106712405582
$ws.Range("C127").Value = 1002.388731936373
$ws.Range("C128").Value = 1128.996380299766
$ws.Range("C129").Value = 2361.056581219794
$ws.Range("C130").Value = 1441.783971398429
$ws.Range("C131").Value = 956.659691840205
$ws.Range("C132").Value = 1463.71052702022
$ws.Range("C133").Value = 5330.539154475424
$ws.Range("C134").Value = 3314.741082534716
$ws.Range("C135").Value = 1625.905825842452
$ws.Range("C136").Value = 3210.869677115934
$ws.Range("C137").Value = 2067.29003376698
$ws.Range("C138").Value = 1751.664428859304
$ws.Range("C139").Value = 1093.134170274031
$ws.Range("C140").Value = 6487.899081675427
$ws.Range("C141").Value = 571.453129531788
$ws.Range("C142").Value = 469.9423670895969
$ws.Range("C143").Value = 579.0880693780265
$ws.Range("C144").Value = 731.9993357350996
$ws.Range("C145").Value = 4944.191641077407
$ws.Range("C146").Value = 6114.227214287786
$ws.Range("C147").Value = 0
$ws.Range("C148").Value = 7556.788578822353
$ws.Range("C149").Value = 871.998368594318
$ws.Range("C150").Value = 1641.006984799246
$ws.Range("C151").Value = 1062.040157863007
$ws.Range("C152").Value = 1134.924536209078
$ws.Range("C153").Value = 2425.561644739583
$ws.Range("C154").Value = 1469.192636109792
$ws.Range("C155").Value = 869.0586852798759
$ws.Range("C156").Value = 1529.507453727912
$ws.Range("C157").Value = 5176.058803160127
$ws.Range("C158").Value = 3382.563653843273
$ws.Range("C159").Value = 1644.598009122967
$ws.Range("C160").Value = 3242.636921959078
$ws.Range("C161").Value = 729.8559996981501
$ws.Range("C162").Value = 2111.193164269742
$ws.Range("C163").Value = 1129.713195979213
$ws.Range("C164").Value = 6411.986543373589
$ws.Range("C165").Value = 548.2681436079887
$ws.Range("C166").Value = 1895.214690888655
$ws.Range("C167").Value = 5089.61202008711
$ws.Range("C168").Value = 6262.368904654469
$ws.Range("C169").Value = 359.6000402964525
$ws.Range("C170").Value = 10236.89594140222
$ws.Range("C171").Value = 0
$ws.Range("C172").Value = 904.381892324943

# --- Update AL flags ---
$ws.Range("AL28").Value = 1
$ws.Range("AL48").Value = 1
$ws.Range("AL161").Value = 1

Write-Host "done"